# Update NATMI TPM ligand/receptor-derived statistics (columns G:T, rows 2:26)
# on Sheet1 with newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4.269024333333333, 12.807073, 0.239419143716247, 0.2624375843396372, 3, 1, 2.906846333333333, 8.720538999999999, 0.005520525738044089, 0.005624540846623205, 12.40939773026077, 111.684579572347, 0.001321719545066018, 0.001476090912807412),
    @(4.269024333333333, 12.807073, 0.239419143716247, 0.2624375843396372, 3, 1, 185.8027443333333, 557.408233, 0.3528665483720876, 0.3595150912979765, 793.1964367591121, 7138.767930832009, 0.08448300685735287, 0.09435027209388509),
    @(4.269024333333333, 12.807073, 0.239419143716247, 0.2624375843396372, 3, 1, 137.0717086666666, 411.2151259999999, 0.2603191943704447, 0.2652240042658267, 585.1624597095774, 5266.462137386197, 0.06232539860907512, 0.0696047469884092),
    @(4.269024333333333, 12.807073, 0.239419143716247, 0.2624375843396372, 2, 1, 29.2127365, 58.425473, 0.05547925319534149, 0.03768304451958546, 124.7098829617548, 748.259297770529, 0.01328279529404552, 0.009889447174283012),
    @(4.269024333333333, 12.807073, 0.239419143716247, 0.2624375843396372, 3, 1, 171.5584106666666, 514.6752319999999, 0.3258144783240821, 0.331953319069988, 732.3870297239928, 6591.483267515936, 0.07800622341070745, 0.08711702717025249),
    @(2.66687, 8.00061, 0.1495657278917394, 0.1639454043592587, 3, 1, 2.906846333333333, 8.720538999999999, 0.005520525738044089, 0.005624540846623205, 7.752181280976664, 69.76963152878999, 0.0008256814503556461, 0.0009221176234348088),
    @(2.66687, 8.00061, 0.1495657278917394, 0.1639454043592587, 3, 1, 185.8027443333333, 557.408233, 0.3528665483720876, 0.3595150912979765, 495.5117647802366, 4459.60588302213, 0.05277674215591696, 0.05894084701610258),
    @(2.66687, 8.00061, 0.1495657278917394, 0.1639454043592587, 3, 1, 137.0717086666666, 411.2151259999999, 0.2603191943704447, 0.2652240042658267, 365.5524276918732, 3289.971849226859, 0.03893482979020674, 0.04348225662514272),
    @(2.66687, 8.00061, 0.1495657278917394, 0.1639454043592587, 2, 1, 29.2127365, 58.425473, 0.05547925319534149, 0.03768304451958546, 77.90657058975499, 467.43942353853, 0.008297794887051358, 0.006177961971251386),
    @(2.66687, 8.00061, 0.1495657278917394, 0.1639454043592587, 3, 1, 171.5584106666666, 514.6752319999999, 0.3258144783240821, 0.331953319069988, 457.5239786546132, 4117.71580789152, 0.0487306796082087, 0.05442222112332722),
    @(3.607962333333333, 10.823887, 0.2023448884238746, 0.221798904202795, 3, 1, 2.906846333333333, 8.720538999999999, 0.005520525738044089, 0.005624540846623205, 10.48779207945477, 94.39012871509298, 0.001117050164505659, 0.001247516996424888),
    @(3.607962333333333, 10.823887, 0.2023448884238746, 0.221798904202795, 3, 1, 185.8027443333333, 557.408233, 0.3528665483720876, 0.3595150912979765, 670.3693029846301, 6033.32372686167, 0.07140074235886783, 0.079740053294259),
    @(3.607962333333333, 10.823887, 0.2023448884238746, 0.221798904202795, 3, 1, 137.0717086666666, 411.2151259999999, 0.2603191943704447, 0.2652240042658267, 494.5495618349734, 4450.946056514761, 0.05267425833948055, 0.05882639351443779),
    @(3.607962333333333, 10.823887, 0.2023448884238746, 0.221798904202795, 2, 1, 29.2127365, 58.425473, 0.05547925319534149, 0.03768304451958546, 105.3984529455918, 632.390717673551, 0.01122594329765126, 0.008358057981469194),
    @(3.607962333333333, 10.823887, 0.2023448884238746, 0.221798904202795, 3, 1, 171.5584106666666, 514.6752319999999, 0.3258144783240821, 0.331953319069988, 618.9762836518647, 5570.786552866783, 0.06592689426336931, 0.07362688241620412),
    @(4.691815, 9.38363, 0.2631304584056895, 0.1922857150526861, 3, 1, 2.906846333333333, 8.720538999999999, 0.005520525738044089, 0.005624540846623205, 13.63838522942833, 81.83031137656999, 0.001452618468091949, 0.001081518858535983),
    @(4.691815, 9.38363, 0.2631304584056895, 0.1922857150526861, 3, 1, 185.8027443333333, 557.408233, 0.3528665483720876, 0.3595150912979765, 871.7521029042983, 5230.51261742579, 0.09284993662918084, 0.06912961640246315),
    @(4.691815, 9.38363, 0.2631304584056895, 0.1922857150526861, 3, 1, 137.0717086666666, 411.2151259999999, 0.2603191943704447, 0.2652240042658267, 643.1150987978965, 3858.69059278738, 0.06849790894649489, 0.05099878730939116),
    @(4.691815, 9.38363, 0.2631304584056895, 0.1922857150526861, 2, 1, 29.2127365, 58.425473, 0.05547925319534149, 0.03768304451958546, 137.0607553017475, 548.24302120699, 0.01459828132529552, 0.007245911160810694),
    @(4.691815, 9.38363, 0.2631304584056895, 0.1922857150526861, 3, 1, 171.5584106666666, 514.6752319999999, 0.3258144783240821, 0.331953319069988, 804.9203245420266, 4829.521947252159, 0.08573171303662633, 0.06382988132148511),
    @(2.595084333333333, 7.785253, 0.1455397815624493, 0.159532392045623, 3, 1, 2.906846333333333, 8.720538999999999, 0.005520525738044089, 0.005624540846623205, 7.543511379040776, 67.89160241136699, 0.0008034561100248161, 0.0008972964554201136),
    @(2.595084333333333, 7.785253, 0.1455397815624493, 0.159532392045623, 3, 1, 185.8027443333333, 557.408233, 0.3528665483720876, 0.3595150912979765, 482.1737909097721, 4339.564118187949, 0.0513561203707691, 0.05735430249126675),
    @(2.595084333333333, 7.785253, 0.1455397815624493, 0.159532392045623, 3, 1, 137.0717086666666, 411.2151259999999, 0.2603191943704447, 0.2652240042658267, 355.7126437040975, 3201.413793336877, 0.03788679868518731, 0.04231181982844586),
    @(2.595084333333333, 7.785253, 0.1455397815624493, 0.159532392045623, 2, 1, 29.2127365, 58.425473, 0.05547925319534149, 0.03768304451958546, 75.80951482494483, 454.857088949669, 0.00807443839129782, 0.006011666231771173),
    @(2.595084333333333, 7.785253, 0.1455397815624493, 0.159532392045623, 3, 1, 171.5584106666666, 514.6752319999999, 0.3258144783240821, 0.331953319069988, 445.2085437726328, 4006.876893953695, 0.0474189680051703, 0.05295730703871912)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 7  # column G is 7
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
